$d = $word.ActiveDocument
$s = $d.Styles("Normal")
try {
  $s.Delete()
  Write-Output "deleted"
} catch {
  Write-Output "delete err: $_"
}
